$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on A3/B3 so numeric-looking strings keep leading zeros
# and aren't converted to actual numbers (matches the other rows' inline string style).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"

$ws.Range("A3").Value = "035148"
$ws.Range("B3").Value = "9014224"
$ws.Range("C3").Value = "001AE87EB516"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "TELEFONO"
$ws.Range("H3").Value = "JORGE ANDRES MELO MAYORGA"
